# Update countries & provincias Spain
# Applies the refreshed COVID-19 snapshot (13 de Mayo de 2020, 05:05) to the
# "Pais" worksheet: a handful of countries changed case counts and therefore
# swapped ranking positions (rows are sorted by total cases, column B, desc),
# plus the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Honduras (row 79): values refreshed in place, no re-ranking needed ---
$ws.Cells.Item(79, 4).Value = 211   # D79 Casos activos
$ws.Cells.Item(79, 5).Value = 1748  # E79 Recuperados
$ws.Cells.Item(79, 7).Value = 5     # G79 Muertes hoy
$ws.Cells.Item(79, 8).Value = 121   # H79 Muertes

# --- Haiti overtakes Nepal and Liberia (rows 142-144) ---
# Row 142 now shows Haiti with its newly updated totals.
$ws.Cells.Item(142, 1).Value = "Haiti"
$ws.Cells.Item(142, 2).Value = 219
$ws.Cells.Item(142, 3).Value = 10
$ws.Cells.Item(142, 4).Value = 17
$ws.Cells.Item(142, 5).Value = 184
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 2
$ws.Cells.Item(142, 8).Value = 18

# Row 143 now shows Nepal (its own figures, unchanged).
$ws.Cells.Item(143, 1).Value = "Nepal"
$ws.Cells.Item(143, 2).Value = 217
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 33
$ws.Cells.Item(143, 5).Value = 184
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

# Row 144 now shows Liberia (its own figures, unchanged).
$ws.Cells.Item(144, 1).Value = "Liberia"
$ws.Cells.Item(144, 2).Value = 211
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 85
$ws.Cells.Item(144, 5).Value = 106
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 20

# --- Belice overtakes Nueva Caledonia (rows 193-194, tied totals) ---
$ws.Cells.Item(193, 1).Value = "Belice"
$ws.Cells.Item(193, 2).Value = 18
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 16
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 2

$ws.Cells.Item(194, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(194, 2).Value = 18
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 18
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 0

# --- Curazao overtakes Dominica (rows 198-199, tied totals) ---
$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 2).Value = 16
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 5).Value = 1
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 1

$ws.Cells.Item(199, 1).Value = "Dominica"
$ws.Cells.Item(199, 2).Value = 16
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 15
$ws.Cells.Item(199, 5).Value = 1
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# --- Sahara Occidental overtakes San Bartolome (rows 215-216, identical totals) ---
$ws.Cells.Item(215, 1).Value = "Sahara Occidental"
$ws.Cells.Item(215, 2).Value = 6
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 6
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

$ws.Cells.Item(216, 1).Value = "San Bartolome"
$ws.Cells.Item(216, 2).Value = 6
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 6
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0

# --- Refresh the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 05:05"
